$wb = $excel.ActiveWorkbook

# --- Existing sheets ---
$ws1 = $wb.Worksheets.Item("Test 1")
$ws2 = $wb.Worksheets.Item("Test 2")

# Sheet "Test 1": selection changes to the whole used range A1:B5
$ws1.Activate()
[void]$ws1.Range("A1:B5").Select()

# Sheet "Test 2": selection changes to A6:B6 (active cell A6)
$ws2.Activate()
[void]$ws2.Range("A6:B6").Select()

# --- New sheet "Test 3", added after "Test 2" ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "Test 3"

# Raw timestamp data (column A) and elapsed-time formulas (column B)
$ws3.Range("A1").Value = 1537059176163
$ws3.Range("A2").Value = 1537066948183
$ws3.Range("A3").Value = 1537075466130
$ws3.Range("A4").Value = 1537083807659

$ws3.Range("B2").Formula = "=(A2-A1)/1000/60/60"
$ws3.Range("B3").Formula = "=(A3-A2)/1000/60/60"
$ws3.Range("B4").Formula = "=(A4-A3)/1000/60/60"

# Row 5: average row, labeled with a single space (matches source workbook)
$ws3.Range("A5").Value = " "
$ws3.Range("B5").Formula = "=AVERAGE(B2:B4)"

# Row 6: delta row, labeled "Delta", formatted as a percentage
$ws3.Range("A6").Value = "Delta"
$ws3.Range("B6").Formula = "=('Test 2'!B5-'Test 3'!B5)/'Test 1'!B5"
$ws3.Range("B6").NumberFormat = "0.00%"

# Selection on the new sheet lands one row below the data (B7), and it
# becomes the active/visible tab.
[void]$ws3.Range("B7").Select()
$ws3.Activate()
